$d = $word.ActiveDocument
$p = $d.Paragraphs(7)
$base = $p.Range.Start
$r1 = $d.Range($base + 70, $base + 71)
Write-Host "before: [" $r1.Text "]"
$r1.Text = "Z"
Write-Host "FINAL:" $d.Paragraphs(7).Range.Text
